$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A24 holds a date-like string ("08/08/2025") that must stay as literal text,
# matching the rest of the "Fecha" column (which stores dates as plain
# strings, not Excel date serials). Assigning the literal directly gets
# auto-parsed into a date serial by Excel's type inference, so we briefly
# force a text number format, assign the value, then clear the format again
# so no stray style is left behind on the cell (matching the target, which
# has no style applied to A24).
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "08/08/2025"
$ws.Range("A24").ClearFormats()

$ws.Range("B24").Value = "U. Magdalena"
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = "Pasto"
$ws.Range("F24").Value = "D"
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 1
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0.68
$ws.Range("L24").Value = 0.73
$ws.Range("M24").Value = 8
$ws.Range("N24").Value = 12
$ws.Range("O24").Value = 2
$ws.Range("P24").Value = 7
